$d = $word.ActiveDocument

# Insert two new paragraphs at the very start of the document:
#   1) "Curso de HTML5" (plain, unformatted)
#   2) an empty paragraph
# Using raw OOXML insertion keeps them free of any inherited paragraph/run
# formatting (e.g. the underline carried by the original paragraph mark).
$startRange = $d.Range(0, 0)
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParagraphsXml = "<w:p $xmlNs><w:r><w:t>Curso de HTML5</w:t></w:r></w:p><w:p $xmlNs/>"
$startRange.InsertXML($newParagraphsXml)

# The original paragraph (with the underline formatting and the _GoBack
# bookmark) is now the 3rd paragraph. Remove its text run but keep the
# paragraph mark formatting and the bookmark start/end intact.
$originalPara = $d.Paragraphs(3)
$originalRange = $d.Range($originalPara.Range.Start, $originalPara.Range.End)
$originalRange.Text = ""
